$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6251971.5
$ws.Range("J40").Value = 38463332
$ws.Range("L40").Value = 38463332
$ws.Range("N40").Value = -38463682

$ws.Range("H132").Value = 1457.5862
$ws.Range("I132").Value = 1217
$ws.Range("J132").Value = 8194
$ws.Range("K132").Value = 3651
$ws.Range("L132").Value = 24582
$ws.Range("M132").Value = -1121
$ws.Range("N132").Value = -29642

$ws.Range("H139").Value = 32460
$ws.Range("J139").Value = 38690
$ws.Range("L139").Value = 38690
$ws.Range("N139").Value = -48970

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6704.653
$ws.Range("I32").Value = 5830.289
$ws.Range("J32").Value = 12752.333
$ws.Range("K32").Value = 5830.289
$ws.Range("L32").Value = 12752.333
$ws.Range("M32").Value = -5543.289
$ws.Range("N32").Value = -13326.333

$ws.Range("H61").Value = 350972.9
$ws.Range("I61").Value = 8484
$ws.Range("J61").Value = 717925.3
$ws.Range("K61").Value = 8484
$ws.Range("L61").Value = 717925.3
$ws.Range("M61").Value = -8272
$ws.Range("N61").Value = -718349.3

$ws.Range("H74").Value = 1740.9429
$ws.Range("I74").Value = 1426.875
$ws.Range("J74").Value = 2426.182
$ws.Range("K74").Value = 1426.875
$ws.Range("L74").Value = 2426.182
$ws.Range("M74").Value = -552.875
$ws.Range("N74").Value = -4174.182

$ws.Range("H77").Value = 1740.9429
$ws.Range("I77").Value = 1426.875
$ws.Range("J77").Value = 2426.182
$ws.Range("K77").Value = 7134.375
$ws.Range("L77").Value = 12130.91
$ws.Range("M77").Value = -2766.375
$ws.Range("N77").Value = -20866.91

$ws.Range("H136").Value = 350972.9
$ws.Range("I136").Value = 8484
$ws.Range("J136").Value = 717925.3
$ws.Range("K136").Value = 25452
$ws.Range("L136").Value = 2153775.9
$ws.Range("M136").Value = -22902
$ws.Range("N136").Value = -2158875.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2070.4614
$ws.Range("I86").Value = 1922.75
$ws.Range("J86").Value = 2306.8
$ws.Range("K86").Value = 1922.75
$ws.Range("L86").Value = 2306.8
$ws.Range("M86").Value = -799.75
$ws.Range("N86").Value = -4552.8

$ws.Range("H89").Value = 2070.4614
$ws.Range("I89").Value = 1922.75
$ws.Range("J89").Value = 2306.8
$ws.Range("K89").Value = 9613.75
$ws.Range("L89").Value = 11534
$ws.Range("M89").Value = -3997.75
$ws.Range("N89").Value = -22766

$ws.Range("H94").Value = 1656.6
$ws.Range("I94").Value = 1146
$ws.Range("J94").Value = 2422.5
$ws.Range("K94").Value = 1146
$ws.Range("L94").Value = 2422.5
$ws.Range("M94").Value = -695
$ws.Range("N94").Value = -3324.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2874.45
$ws.Range("I132").Value = 2486
$ws.Range("J132").Value = 3192.2727
$ws.Range("K132").Value = 7458
$ws.Range("L132").Value = 9576.8181
$ws.Range("M132").Value = -4928
$ws.Range("N132").Value = -14636.8181

$ws.Range("H134").Value = 246771.36
$ws.Range("I134").Value = 2836.0557
$ws.Range("K134").Value = 8508.167099999999
$ws.Range("M134").Value = -5973.167099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1380.1316
$ws.Range("I5").Value = 421.75
$ws.Range("J5").Value = 1822.4615
$ws.Range("K5").Value = 1265.25
$ws.Range("L5").Value = 5467.3845
$ws.Range("M5").Value = -1153.25
$ws.Range("N5").Value = -5691.3845

$ws.Range("H7").Value = 1502.75
$ws.Range("I7").Value = 111
$ws.Range("J7").Value = 1966.6666
$ws.Range("K7").Value = 333
$ws.Range("L7").Value = 5899.9998
$ws.Range("M7").Value = -221
$ws.Range("N7").Value = -6123.9998

$ws.Range("H52").Value = 975
$ws.Range("J52").Value = 975
$ws.Range("L52").Value = 2925
$ws.Range("N52").Value = -3457

$ws.Range("H109").Value = 1935.5834
$ws.Range("I109").Value = 242.33333
$ws.Range("K109").Value = 726.99999
$ws.Range("M109").Value = 313.00001

$ws.Range("H113").Value = 1396002.1
$ws.Range("I113").Value = 1316451.8
$ws.Range("J113").Value = 2000585.2
$ws.Range("K113").Value = 3949355.4
$ws.Range("L113").Value = 6001755.6
$ws.Range("M113").Value = -3947185.4
$ws.Range("N113").Value = -6006095.6

$ws.Range("H122").Value = 3088.814
$ws.Range("I122").Value = 420.5
$ws.Range("K122").Value = 3784.5
$ws.Range("M122").Value = -1334.5

$ws.Range("H131").Value = 2128528.8
$ws.Range("I131").Value = 5000486.5
$ws.Range("K131").Value = 15001459.5
$ws.Range("M131").Value = -14996419.5

$ws.Range("H135").Value = 1380.1316
$ws.Range("I135").Value = 421.75
$ws.Range("J135").Value = 1822.4615
$ws.Range("K135").Value = 3795.75
$ws.Range("L135").Value = 16402.1535
$ws.Range("M135").Value = -1260.75
$ws.Range("N135").Value = -21472.1535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5039.6665
$ws.Range("I70").Value = 5099.909
$ws.Range("K70").Value = 5099.909
$ws.Range("M70").Value = -4829.909

$ws.Range("H73").Value = 5039.6665
$ws.Range("I73").Value = 5099.909
$ws.Range("K73").Value = 5099.909
$ws.Range("M73").Value = -4163.909

$ws.Range("H80").Value = 7364.25
$ws.Range("I80").Value = 9983.462
$ws.Range("K80").Value = 9983.462
$ws.Range("M80").Value = -8985.462

$ws.Range("H83").Value = 7364.25
$ws.Range("I83").Value = 9983.462
$ws.Range("K83").Value = 49917.31
$ws.Range("M83").Value = -44925.31

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 52634460
$ws.Range("I40").Value = 62502930
$ws.Range("K40").Value = 62502930
$ws.Range("M40").Value = -62502794

$ws.Range("H122").Value = 2147472.5
$ws.Range("I122").Value = 2752684
$ws.Range("J122").Value = 836181.2
$ws.Range("K122").Value = 8258052
$ws.Range("L122").Value = 2508543.6
$ws.Range("M122").Value = -8255602
$ws.Range("N122").Value = -2513443.6

$ws.Range("H132").Value = 13896390
$ws.Range("I132").Value = 15880941
$ws.Range("J132").Value = 4533
$ws.Range("K132").Value = 47642823
$ws.Range("L132").Value = 13599
$ws.Range("M132").Value = -47640293
$ws.Range("N132").Value = -18659

$ws.Range("H136").Value = 6464.5093
$ws.Range("I136").Value = 4690.875
$ws.Range("J136").Value = 11921.846
$ws.Range("K136").Value = 14072.625
$ws.Range("L136").Value = 35765.538
$ws.Range("M136").Value = -11522.625
$ws.Range("N136").Value = -40865.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 72099.14
$ws.Range("I62").Value = 83157.336
$ws.Range("J62").Value = 5750
$ws.Range("K62").Value = 83157.336
$ws.Range("L62").Value = 5750
$ws.Range("M62").Value = -82533.336
$ws.Range("N62").Value = -6998

$ws.Range("H65").Value = 72099.14
$ws.Range("I65").Value = 83157.336
$ws.Range("J65").Value = 5750
$ws.Range("K65").Value = 415786.68
$ws.Range("L65").Value = 28750
$ws.Range("M65").Value = -412666.68
$ws.Range("N65").Value = -34990

$ws.Range("H107").Value = 62500732
$ws.Range("I107").Value = 83333980
$ws.Range("J107").Value = 974
$ws.Range("K107").Value = 250001940
$ws.Range("L107").Value = 2922
$ws.Range("M107").Value = -250000020
$ws.Range("N107").Value = -6762

$ws.Range("H122").Value = 1684.1666
$ws.Range("I122").Value = 1381.0588
$ws.Range("J122").Value = 2080.5386
$ws.Range("K122").Value = 4143.1764
$ws.Range("L122").Value = 6241.6158
$ws.Range("M122").Value = -1693.1764
$ws.Range("N122").Value = -11141.6158

$ws.Range("H126").Value = 694.72
$ws.Range("I126").Value = 578.95
$ws.Range("J126").Value = 1157.8
$ws.Range("K126").Value = 1736.85
$ws.Range("L126").Value = 3473.4
$ws.Range("M126").Value = 733.1499999999999
$ws.Range("N126").Value = -8413.4

$ws.Range("H132").Value = 1981.5
$ws.Range("I132").Value = 1405.5
$ws.Range("J132").Value = 2941.5
$ws.Range("K132").Value = 4216.5
$ws.Range("L132").Value = 8824.5
$ws.Range("M132").Value = -1686.5
$ws.Range("N132").Value = -13884.5
